$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 7001.6113
$ws.Range("I76").Value = 8559
$ws.Range("J76").Value = 5054.875
$ws.Range("K76").Value = 8559
$ws.Range("L76").Value = 5054.875
$ws.Range("M76").Value = -8244
$ws.Range("N76").Value = -5684.875
$ws.Range("H79").Value = 7001.6113
$ws.Range("I79").Value = 8559
$ws.Range("J79").Value = 5054.875
$ws.Range("K79").Value = 8559
$ws.Range("L79").Value = 5054.875
$ws.Range("M79").Value = -7467
$ws.Range("N79").Value = -7238.875
$ws.Range("H86").Value = 3614.25
$ws.Range("I86").Value = 1770.5217
$ws.Range("J86").Value = 6876.231
$ws.Range("K86").Value = 1770.5217
$ws.Range("L86").Value = 6876.231
$ws.Range("M86").Value = -647.5217
$ws.Range("N86").Value = -9122.231
$ws.Range("H89").Value = 3614.25
$ws.Range("I89").Value = 1770.5217
$ws.Range("J89").Value = 6876.231
$ws.Range("K89").Value = 8852.6085
$ws.Range("L89").Value = 34381.155
$ws.Range("M89").Value = -3236.6085
$ws.Range("N89").Value = -45613.155
$ws.Range("H98").Value = 894.64
$ws.Range("I98").Value = 582.6842
$ws.Range("J98").Value = 1882.5
$ws.Range("K98").Value = 582.6842
$ws.Range("L98").Value = 1882.5
$ws.Range("M98").Value = 915.3158
$ws.Range("N98").Value = -4878.5
$ws.Range("H106").Value = 3746.1667
$ws.Range("I106").Value = 1790.8
$ws.Range("K106").Value = 1790.8
$ws.Range("M106").Value = -1159.8
$ws.Range("H122").Value = 894.64
$ws.Range("I122").Value = 582.6842
$ws.Range("J122").Value = 1882.5
$ws.Range("K122").Value = 1748.0526
$ws.Range("L122").Value = 5647.5
$ws.Range("M122").Value = 701.9474
$ws.Range("N122").Value = -10547.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7417.2676
$ws.Range("I32").Value = 4693.056
$ws.Range("J32").Value = 20311.867
$ws.Range("K32").Value = 4693.056
$ws.Range("L32").Value = 20311.867
$ws.Range("M32").Value = -4406.056
$ws.Range("N32").Value = -20885.867
$ws.Range("H74").Value = 2189.6592
$ws.Range("I74").Value = 1277.5161
$ws.Range("J74").Value = 4364.769
$ws.Range("K74").Value = 1277.5161
$ws.Range("L74").Value = 4364.769
$ws.Range("M74").Value = -403.5161000000001
$ws.Range("N74").Value = -6112.769
$ws.Range("H77").Value = 2189.6592
$ws.Range("I77").Value = 1277.5161
$ws.Range("J77").Value = 4364.769
$ws.Range("K77").Value = 6387.5805
$ws.Range("L77").Value = 21823.845
$ws.Range("M77").Value = -2019.5805
$ws.Range("N77").Value = -30559.845
$ws.Range("H88").Value = 1925.46
$ws.Range("I88").Value = 1901.825
$ws.Range("J88").Value = 2020
$ws.Range("K88").Value = 1901.825
$ws.Range("L88").Value = 2020
$ws.Range("M88").Value = -1495.825
$ws.Range("N88").Value = -2832
$ws.Range("H91").Value = 1925.46
$ws.Range("I91").Value = 1901.825
$ws.Range("J91").Value = 2020
$ws.Range("K91").Value = 1901.825
$ws.Range("L91").Value = 2020
$ws.Range("M91").Value = -497.825
$ws.Range("N91").Value = -4828
$ws.Range("H132").Value = 3162.5625
$ws.Range("I132").Value = 3087.4517
$ws.Range("J132").Value = 3299.5293
$ws.Range("K132").Value = 9262.355100000001
$ws.Range("L132").Value = 9898.5879
$ws.Range("M132").Value = -6732.355100000001
$ws.Range("N132").Value = -14958.5879

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8109.75
$ws.Range("I86").Value = 7289.5557
$ws.Range("J86").Value = 9164.286
$ws.Range("K86").Value = 7289.5557
$ws.Range("L86").Value = 9164.286
$ws.Range("M86").Value = -6166.5557
$ws.Range("N86").Value = -11410.286
$ws.Range("H89").Value = 8109.75
$ws.Range("I89").Value = 7289.5557
$ws.Range("J89").Value = 9164.286
$ws.Range("K89").Value = 36447.7785
$ws.Range("L89").Value = 45821.43
$ws.Range("M89").Value = -30831.7785
$ws.Range("N89").Value = -57053.43
$ws.Range("H105").Value = 1886.037
$ws.Range("I105").Value = 1760.625
$ws.Range("J105").Value = 2244.3572
$ws.Range("K105").Value = 1760.625
$ws.Range("L105").Value = 2244.3572
$ws.Range("M105").Value = -13.625
$ws.Range("N105").Value = -5738.3572

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2163.4856
$ws.Range("I31").Value = 1308.766
$ws.Range("J31").Value = 3910.087
$ws.Range("K31").Value = 1308.766
$ws.Range("L31").Value = 3910.087
$ws.Range("M31").Value = -1013.766
$ws.Range("N31").Value = -4500.087
$ws.Range("H34").Value = 2163.4856
$ws.Range("I34").Value = 1308.766
$ws.Range("J34").Value = 3910.087
$ws.Range("K34").Value = 1308.766
$ws.Range("L34").Value = 3910.087
$ws.Range("M34").Value = -1106.766
$ws.Range("N34").Value = -4314.087
$ws.Range("H58").Value = 1483.3334
$ws.Range("I58").Value = 1019.7241
$ws.Range("J58").Value = 2517.5386
$ws.Range("K58").Value = 1019.7241
$ws.Range("L58").Value = 2517.5386
$ws.Range("M58").Value = -816.7241
$ws.Range("N58").Value = -2923.5386
$ws.Range("H136").Value = 1483.3334
$ws.Range("I136").Value = 1019.7241
$ws.Range("J136").Value = 2517.5386
$ws.Range("K136").Value = 3059.1723
$ws.Range("L136").Value = 7552.6158
$ws.Range("M136").Value = -509.1723000000002
$ws.Range("N136").Value = -12652.6158
$ws.Range("H140").Value = 22780
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 22780
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 22780
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -33140

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1275.0555
$ws.Range("I107").Value = 155.66667
$ws.Range("J107").Value = 1834.75
$ws.Range("K107").Value = 467.00001
$ws.Range("L107").Value = 5504.25
$ws.Range("M107").Value = 1452.99999
$ws.Range("N107").Value = -9344.25
$ws.Range("H131").Value = 890.0893
$ws.Range("J131").Value = 911.6731
$ws.Range("L131").Value = 2735.0193
$ws.Range("N131").Value = -12815.0193
$ws.Range("H137").Value = 3716.5833
$ws.Range("I137").Value = 3750
$ws.Range("J137").Value = 3699.875
$ws.Range("K137").Value = 11250
$ws.Range("L137").Value = 11099.625
$ws.Range("M137").Value = -6150
$ws.Range("N137").Value = -21299.625

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7652.1816
$ws.Range("I70").Value = 4274.9033
$ws.Range("K70").Value = 4274.9033
$ws.Range("M70").Value = -4004.9033
$ws.Range("H73").Value = 7652.1816
$ws.Range("I73").Value = 4274.9033
$ws.Range("K73").Value = 4274.9033
$ws.Range("M73").Value = -3338.9033
$ws.Range("H80").Value = 4604.6875
$ws.Range("I80").Value = 4937.037
$ws.Range("J80").Value = 2810
$ws.Range("K80").Value = 4937.037
$ws.Range("L80").Value = 2810
$ws.Range("M80").Value = -3939.037
$ws.Range("N80").Value = -4806
$ws.Range("H83").Value = 4604.6875
$ws.Range("I83").Value = 4937.037
$ws.Range("J83").Value = 2810
$ws.Range("K83").Value = 24685.185
$ws.Range("L83").Value = 14050
$ws.Range("M83").Value = -19693.185
$ws.Range("N83").Value = -24034
$ws.Range("H97").Value = 2026.1
$ws.Range("I97").Value = 1200
$ws.Range("J97").Value = 2852.2
$ws.Range("K97").Value = 1200
$ws.Range("L97").Value = 2852.2
$ws.Range("M97").Value = -704
$ws.Range("N97").Value = -3844.2
$ws.Range("H122").Value = 1041.8889
$ws.Range("I122").Value = 878.2
$ws.Range("J122").Value = 1246.5
$ws.Range("K122").Value = 2634.6
$ws.Range("L122").Value = 3739.5
$ws.Range("M122").Value = -184.6000000000004
$ws.Range("N122").Value = -8639.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 7787.75
$ws.Range("I53").Value = 6366.6665
$ws.Range("J53").Value = 12051
$ws.Range("K53").Value = 6366.6665
$ws.Range("L53").Value = 12051
$ws.Range("M53").Value = -5848.6665
$ws.Range("N53").Value = -13087
$ws.Range("H93").Value = 1794.2
$ws.Range("I93").Value = 1705.6
$ws.Range("J93").Value = 2060
$ws.Range("K93").Value = 1705.6
$ws.Range("L93").Value = 2060
$ws.Range("M93").Value = -457.5999999999999
$ws.Range("N93").Value = -4556
$ws.Range("H132").Value = 12614.096
$ws.Range("I132").Value = 3900.3
$ws.Range("J132").Value = 20535.727
$ws.Range("K132").Value = 11700.9
$ws.Range("L132").Value = 61607.181
$ws.Range("M132").Value = -9170.900000000001
$ws.Range("N132").Value = -66667.181
$ws.Range("H136").Value = 4707.442
$ws.Range("I136").Value = 2522.8147
$ws.Range("J136").Value = 8394
$ws.Range("K136").Value = 7568.4441
$ws.Range("L136").Value = 25182
$ws.Range("M136").Value = -5018.4441
$ws.Range("N136").Value = -30282

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 20084
$ws.Range("J50").Value = 20084
$ws.Range("L50").Value = 20084
$ws.Range("N50").Value = -21346
$ws.Range("H53").Value = 7999.5
$ws.Range("J53").Value = 7999.5
$ws.Range("L53").Value = 7999.5
$ws.Range("M53").Value = -9213.5
$ws.Range("H96").Value = 10642
$ws.Range("I96").Value = 1549.6666
$ws.Range("J96").Value = 19734.334
$ws.Range("K96").Value = 1549.6666
$ws.Range("L96").Value = 19734.334
$ws.Range("M96").Value = -176.6666
$ws.Range("N96").Value = -22480.334
$ws.Range("H122").Value = 54201.05
$ws.Range("I122").Value = 92276.55
$ws.Range("J122").Value = 1847.25
$ws.Range("K122").Value = 276829.65
$ws.Range("L122").Value = 5541.75
$ws.Range("M122").Value = -274379.65
$ws.Range("N122").Value = -10441.75
$ws.Range("H132").Value = 30570.361
$ws.Range("I132").Value = 127267.75
$ws.Range("J132").Value = 2942.5356
$ws.Range("K132").Value = 381803.25
$ws.Range("L132").Value = 8827.606800000001
$ws.Range("M132").Value = -379273.25
$ws.Range("N132").Value = -13887.6068
$ws.Range("H136").Value = 22600686
$ws.Range("I136").Value = 25001264
$ws.Range("J136").Value = 17546834
$ws.Range("K136").Value = 75003792
$ws.Range("L136").Value = 52640502
$ws.Range("M136").Value = -75001242
$ws.Range("N136").Value = -52645602
